$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update word labels in column B
$ws.Range("B9").Value = "<brog>"
$ws.Range("B13").Value = "<command>"

# Update counts in column C
$ws.Range("C3").Value = 10
$ws.Range("C4").Value = 14
$ws.Range("C5").Value = 14
$ws.Range("C6").Value = 14
$ws.Range("C7").Value = 11
$ws.Range("C8").Value = 15
$ws.Range("C9").Value = 21
$ws.Range("C10").Value = 10
$ws.Range("C11").Value = 11
$ws.Range("C12").Value = 12
$ws.Range("C13").Value = 10
$ws.Range("C14").Value = 14
$ws.Range("C16").Value = 10
$ws.Range("C17").Value = 16
$ws.Range("C18").Value = 10
